$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header values for the new columns (N15:AS15)
$ws.Range("N15").Value = 'barometric_press'
$ws.Range("O15").Value = 'carb_dioxide'
$ws.Range("P15").Value = 'carb_monoxide'
$ws.Range("Q15").Value = 'chem_administration'
$ws.Range("R15").Value = 'elev'
$ws.Range("S15").Value = 'humidity'
$ws.Range("T15").Value = 'isolation_source'
$ws.Range("U15").Value = 'methane'
$ws.Range("V15").Value = 'misc_param'
$ws.Range("W15").Value = 'organism_count'
$ws.Range("X15").Value = 'oxy_stat_samp'
$ws.Range("Y15").Value = 'oxygen'
$ws.Range("Z15").Value = 'perturbation'
$ws.Range("AA15").Value = 'pollutants'
$ws.Range("AB15").Value = 'rel_to_oxygen'
$ws.Range("AC15").Value = 'resp_part_matter'
$ws.Range("AD15").Value = 'samp_collect_device'
$ws.Range("AE15").Value = 'samp_mat_process'
$ws.Range("AF15").Value = 'samp_salinity'
$ws.Range("AG15").Value = 'samp_size'
$ws.Range("AH15").Value = 'samp_store_dur'
$ws.Range("AI15").Value = 'samp_store_loc'
$ws.Range("AJ15").Value = 'samp_store_temp'
$ws.Range("AK15").Value = 'samp_vol_we_dna_ext'
$ws.Range("AL15").Value = 'solar_irradiance'
$ws.Range("AM15").Value = 'source_material_id'
$ws.Range("AN15").Value = 'temperature'
$ws.Range("AO15").Value = 'ventilation_rate'
$ws.Range("AP15").Value = 'ventilation_type'
$ws.Range("AQ15").Value = 'volatile_org_comp'
$ws.Range("AR15").Value = 'wind_direction'
$ws.Range("AS15").Value = 'wind_speed'

# Copy the "optional field" header style (from C15) onto the new header cells
$ws.Range("C15").Copy()
$ws.Range("N15:AS15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the definition comments for each new header cell
$ws.Range("N15").AddComment('force per unit area exerted against a surface by the weight of air above that surface')
$ws.Range("O15").AddComment('carbon dioxide (gas) amount or concentration at the time of sampling')
$ws.Range("P15").AddComment('carbon monoxide (gas) amount or concentration at the time of sampling')
$ws.Range("Q15").AddComment('list of chemical compounds administered to the host or site where sampling occurred, and when (e.g. antibiotics, N fertilizer, air filter); can include multiple compounds. For Chemical Entities of Biological Interest ontology (CHEBI) (v1.72), please see http://bioportal.bioontology.org/visualize/44603')
$ws.Range("R15").AddComment('The elevation of the sampling site as measured by the vertical distance from mean sea level.')
$ws.Range("S15").AddComment('amount of water vapour in the air, at the time of sampling')
$ws.Range("T15").AddComment('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("U15").AddComment('methane (gas) amount or concentration at the time of sampling')
$ws.Range("V15").AddComment('any other measurement performed or parameter collected, that is not listed here')
$ws.Range("W15").AddComment('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$ws.Range("X15").AddComment('oxygenation status of sample')
$ws.Range("Y15").AddComment('oxygen (gas) amount or concentration at the time of sampling')
$ws.Range("Z15").AddComment('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$ws.Range("AA15").AddComment('pollutant types and, amount or concentrations measured at the time of sampling; can report multiple pollutants by entering numeric values preceded by name of pollutant')
$ws.Range("AB15").AddComment('Aerobic or anaerobic')
$ws.Range("AC15").AddComment('concentration of substances that remain suspended in the air, and comprise mixtures of organic and inorganic substances (PM10 and PM2.5); can report multiple PM''s by entering numeric values preceded by name of PM')
$ws.Range("AD15").AddComment('Method or device employed for collecting sample')
$ws.Range("AE15").AddComment('Processing applied to the sample during or after isolation')
$ws.Range("AF15").AddComment('salinity of sample, i.e. measure of total salt concentration')
$ws.Range("AG15").AddComment('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("AH15").AddComment('duration for which sample was stored')
$ws.Range("AI15").AddComment('location at which sample was stored, usually name of a specific freezer/room')
$ws.Range("AJ15").AddComment('temperature at which sample was stored, e.g. -80')
$ws.Range("AK15").AddComment('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("AL15").AddComment('the amount of solar energy that arrives at a specific area of a surface during a specific time interval')
$ws.Range("AM15").AddComment('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("AN15").AddComment('temperature of the sample at time of sampling')
$ws.Range("AO15").AddComment('ventilation rate of the system in the sampled premises')
$ws.Range("AP15").AddComment('ventilation system used in the sampled premises')
$ws.Range("AQ15").AddComment('concentration of carbon-based chemicals that easily evaporate at room temperature; can report multiple volatile organic compounds by entering numeric values preceded by name of compound')
$ws.Range("AR15").AddComment('wind direction is the direction from which a wind originates')
$ws.Range("AS15").AddComment('speed of wind measured at the time of sampling')
